{"js": "// Append the new \"logistic\" model block to the end of the document body,\n// right after the existing trailing empty paragraph (and before sectPr).\n// Matches the diff: 2 blank paragraphs, then the logistic_control line,\n// a blank line, set.seed(100), the logistic <- train(...) line, a blank\n// line, summary(logistic), a blank line, and the confusionMatrix line.\n\nconst body = context.document.body;\n\n// Lines to append, in order. `null` marks a blank paragraph.\nconst lines = [\n  null,\n  null,\n  'logistic_control <- trainControl(method = \"cv\", number = 10, sampling = \"down\"), preProc = c(\"center\",\"scale\")',\n  null,\n  \"set.seed(100)\",\n  'logistic <- train(diabetes ~ ., data = train_data, method = \"glm\", family = \"binomial\", trControl = logistic_control)',\n  null,\n  \"summary(logistic)\",\n  null,\n  \"confusionMatrix(logistic) #Accuracy (average) : 0.7319\",\n];\n\nfor (const line of lines) {\n  body.insertParagraph(line ?? \"\", \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Append the new \"logistic\" model block to the end of the document,\n# right after the existing trailing empty paragraph (and before sectPr).\n# Matches the diff: 2 blank paragraphs, then the logistic_control line,\n# a blank line, set.seed(100), the logistic <- train(...) line, a blank\n# line, summary(logistic), a blank line, and the confusionMatrix line.\n\n$d = $word.ActiveDocument\n\n$lines = @(\n    $null,\n    $null,\n    'logistic_control <- trainControl(method = \"cv\", number = 10, sampling = \"down\"), preProc = c(\"center\",\"scale\")',\n    $null,\n    \"set.seed(100)\",\n    'logistic <- train(diabetes ~ ., data = train_data, method = \"glm\", family = \"binomial\", trControl = logistic_control)',\n    $null,\n    \"summary(logistic)\",\n    $null,\n    \"confusionMatrix(logistic) #Accuracy (average) : 0.7319\"\n)\n\nforeach ($line in $lines) {\n    $r = $d.Paragraphs.Last.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    if ($line) {\n        $newRange = $d.Paragraphs.Last.Range\n        $newRange.Text = $line\n    }\n}\n"}
